$d = $word.ActiveDocument

function Replace-ExactText($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $find"
    }
}

Replace-ExactText "811÷4=202, 3" "742÷8=92, 6"
Replace-ExactText "159÷5=31, 4" "241÷2=120, 1"
Replace-ExactText "426÷2=213, 0" "590÷2=295, 0"
Replace-ExactText "527÷8=65, 7" "343÷3=114, 1"
Replace-ExactText "846÷6=141, 0" "606÷6=101, 0"
Replace-ExactText "914÷9=101, 5" "182÷6=30, 2"
Replace-ExactText "858÷2=429, 0" "628÷2=314, 0"
Replace-ExactText "585÷9=65, 0" "609÷5=121, 4"
Replace-ExactText "803÷5=160, 3" "888÷7=126, 6"
Replace-ExactText "532÷2=266, 0" "975÷9=108, 3"
Replace-ExactText "535÷4=133, 3" "454÷9=50, 4"
Replace-ExactText "163÷5=32, 3" "515÷8=64, 3"
Replace-ExactText "447÷3=149, 0" "417÷6=69, 3"
Replace-ExactText "812÷6=135, 2" "397÷6=66, 1"
Replace-ExactText "653÷9=72, 5" "992÷8=124, 0"
Replace-ExactText "145÷2=72, 1" "102÷4=25, 2"
Replace-ExactText "356÷9=39, 5" "143÷4=35, 3"
Replace-ExactText "710÷3=236, 2" "809÷4=202, 1"
Replace-ExactText "314÷5=62, 4" "798÷8=99, 6"
Replace-ExactText "378÷9=42, 0" "312÷5=62, 2"
Replace-ExactText "349÷9=38, 7" "864÷4=216, 0"
Replace-ExactText "891÷6=148, 3" "135÷9=15, 0"
Replace-ExactText "645÷6=107, 3" "508÷3=169, 1"
Replace-ExactText "934÷9=103, 7" "649÷4=162, 1"
Replace-ExactText "394÷4=98, 2" "312÷2=156, 0"

Write-Output "Done: replaced 25 cells"
